{"js": "// Map of old text -> new text, exactly as they appear (unique, whole-run) strings.\nconst replacements = [\n  [\"2024-09-27 Friday\", \"2024-09-28 Saturday\"],\n  [\"72\u00d792=\", \"40\u00d767=\"],\n  [\"69\u00d733=\", \"77\u00d772=\"],\n  [\"34\u00d734=\", \"86\u00d783=\"],\n  [\"72\u00d770=\", \"64\u00d711=\"],\n  [\"99\u00d739=\", \"76\u00d742=\"],\n  [\"75\u00d734=\", \"22\u00d726=\"],\n  [\"32\u00d720=\", \"94\u00d770=\"],\n  [\"65\u00d770=\", \"70\u00d792=\"],\n  [\"86\u00d769=\", \"57\u00d743=\"],\n  [\"51\u00d729=\", \"56\u00d797=\"],\n  [\"65\u00d749=\", \"80\u00d761=\"],\n  [\"52\u00d769=\", \"35\u00d721=\"],\n  [\"87\u00d735=\", \"19\u00d759=\"],\n  [\"90\u00d749=\", \"63\u00d778=\"],\n  [\"63\u00d742=\", \"47\u00d753=\"],\n  [\"74\u00d746=\", \"39\u00d757=\"],\n  [\"36\u00d734=\", \"53\u00d721=\"],\n  [\"27\u00d735=\", \"56\u00d791=\"],\n  [\"79\u00d752=\", \"43\u00d714=\"],\n  [\"52\u00d770=\", \"53\u00d720=\"],\n  [\"56\u00d782=\", \"12\u00d714=\"],\n  [\"64\u00d715=\", \"39\u00d723=\"],\n  [\"47\u00d743=\", \"50\u00d722=\"],\n  [\"11\u00d719=\", \"27\u00d745=\"],\n  [\"92\u00d723=\", \"69\u00d784=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-09-27 Friday\", \"2024-09-28 Saturday\"),\n  @(\"72\u00d792=\", \"40\u00d767=\"),\n  @(\"69\u00d733=\", \"77\u00d772=\"),\n  @(\"34\u00d734=\", \"86\u00d783=\"),\n  @(\"72\u00d770=\", \"64\u00d711=\"),\n  @(\"99\u00d739=\", \"76\u00d742=\"),\n  @(\"75\u00d734=\", \"22\u00d726=\"),\n  @(\"32\u00d720=\", \"94\u00d770=\"),\n  @(\"65\u00d770=\", \"70\u00d792=\"),\n  @(\"86\u00d769=\", \"57\u00d743=\"),\n  @(\"51\u00d729=\", \"56\u00d797=\"),\n  @(\"65\u00d749=\", \"80\u00d761=\"),\n  @(\"52\u00d769=\", \"35\u00d721=\"),\n  @(\"87\u00d735=\", \"19\u00d759=\"),\n  @(\"90\u00d749=\", \"63\u00d778=\"),\n  @(\"63\u00d742=\", \"47\u00d753=\"),\n  @(\"74\u00d746=\", \"39\u00d757=\"),\n  @(\"36\u00d734=\", \"53\u00d721=\"),\n  @(\"27\u00d735=\", \"56\u00d791=\"),\n  @(\"79\u00d752=\", \"43\u00d714=\"),\n  @(\"52\u00d770=\", \"53\u00d720=\"),\n  @(\"56\u00d782=\", \"12\u00d714=\"),\n  @(\"64\u00d715=\", \"39\u00d723=\"),\n  @(\"47\u00d743=\", \"50\u00d722=\"),\n  @(\"11\u00d719=\", \"27\u00d745=\"),\n  @(\"92\u00d723=\", \"69\u00d784=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
